# Corrected CAT error rates report
# Missing "errcord" column caused rates to show as Inf/NaN instead of
# the correct computed error-rate values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "5.5 (1.2, 1.4)"
$ws.Range("D3").Value = "0 (0, 0)"
$ws.Range("D4").Value = "6 (1.2, 1.8)"
$ws.Range("D5").Value = "0.6 (0.3, 0.3)"
